# regen sval data to filter save games
# Update the numeric stat columns (B:G) for rows 2-5 on the single data sheet
# with the newly regenerated values, leaving labels/headers untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.2917716402565462
$ws.Range("C2").Value = 0.04071648406533734
$ws.Range("D2").Value = 6708.013860684405
$ws.Range("E2").Value = 1133.036916526867
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 7841.383265335594

$ws.Range("B3").Value = 0.04271373187048222
$ws.Range("C3").Value = 250555.8564151394
$ws.Range("D3").Value = 3.537761648806719
$ws.Range("E3").Value = 1133.036916526867
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 251692.4738070469

$ws.Range("B4").Value = 1.455362044514542
$ws.Range("C4").Value = 1.655778082260271
$ws.Range("D4").Value = 22.3905356188092
$ws.Range("E4").Value = 10.19245300693656
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 35.69412875252057

$ws.Range("B5").Value = 1.455362044514542
$ws.Range("C5").Value = 0.306821227259698
$ws.Range("D5").Value = 0.1494219747398047
$ws.Range("E5").Value = 0.4942365360607697
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 2.405841782574814
